$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" column (C) for rows 2-18 from 45224 to 45233,
# keeping the existing date formatting/style untouched.
$ws.Range("C2:C18").Value = 45233
